# Apply the edit described by the diff:
# - Remove the erroneous "CROANCA" row (old row 3), shifting the other
#   topic/social rows up by one.
# - Update the negativo/positivo numeric values to their new totals
#   (the dataset behind the rows was recalculated / rows were added
#   upstream, per the commit message "added rows to dataset").
# - The last row (old row 12) disappears, shrinking the used range
#   from A1:E12 to A1:E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 3 ("CROANCA") entirely; this shifts rows 4-12 up to 3-11.
$ws.Rows.Item(3).Delete()

# Now rewrite B:E for rows 3-11 with the final values from the updated dataset.
$ws.Range("B3").Value = "CRONACA"
$ws.Range("C3").Value = "Facebook"
$ws.Range("D3").Value = 840
$ws.Range("E3").Value = 356

$ws.Range("B4").Value = "CRONACA"
$ws.Range("C4").Value = "Instagram"
$ws.Range("D4").Value = 1007
$ws.Range("E4").Value = 193

$ws.Range("B5").Value = "CRONACA"
$ws.Range("C5").Value = "YouTube"
$ws.Range("D5").Value = 913
$ws.Range("E5").Value = 276

$ws.Range("B6").Value = "CRONACA NERA"
$ws.Range("C6").Value = "Facebook"
$ws.Range("D6").Value = 1064
$ws.Range("E6").Value = 136

$ws.Range("B7").Value = "CRONACA NERA"
$ws.Range("C7").Value = "Instagram"
$ws.Range("D7").Value = 1044
$ws.Range("E7").Value = 156

$ws.Range("B8").Value = "CRONACA NERA"
$ws.Range("C8").Value = "YouTube"
$ws.Range("D8").Value = 1023
$ws.Range("E8").Value = 177

$ws.Range("B9").Value = "POLITICA"
$ws.Range("C9").Value = "Facebook"
$ws.Range("D9").Value = 874
$ws.Range("E9").Value = 325

$ws.Range("B10").Value = "POLITICA"
$ws.Range("C10").Value = "Instagram"
$ws.Range("D10").Value = 974
$ws.Range("E10").Value = 226

$ws.Range("B11").Value = "POLITICA"
$ws.Range("C11").Value = "YouTube"
$ws.Range("D11").Value = 895
$ws.Range("E11").Value = 300
